# Applies numeric corrections to the Leve profit-tracking sheets
# (scheduled-runner refresh of market-board prices/profits).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2566.6667
$ws.Range("I64").Value = 2566.6667
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2566.6667
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -2318.6667
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 2566.6667
$ws.Range("I67").Value = 2566.6667
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2566.6667
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -1708.6667
$ws.Range("N67").ClearContents()

$ws.Range("H132").Value = 1860.5385
$ws.Range("I132").Value = 1990.1818
$ws.Range("J132").Value = 1147.5
$ws.Range("K132").Value = 5970.5454
$ws.Range("L132").Value = 3442.5
$ws.Range("M132").Value = -3440.5454
$ws.Range("N132").Value = -8502.5

$ws.Range("H134").Value = 44520
$ws.Range("J134").Value = 44520
$ws.Range("L134").Value = 44520
$ws.Range("N134").Value = -54660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26030.33
$ws.Range("I32").Value = 24614.629
$ws.Range("J32").Value = 31193.47
$ws.Range("K32").Value = 24614.629
$ws.Range("L32").Value = 31193.47
$ws.Range("M32").Value = -24327.629
$ws.Range("N32").Value = -31767.47

$ws.Range("H98").Value = 21087.5
$ws.Range("J98").Value = 21087.5
$ws.Range("L98").Value = 21087.5
$ws.Range("N98").Value = -27077.5

$ws.Range("H134").Value = 58032.25
$ws.Range("J134").Value = 58032.25
$ws.Range("L134").Value = 58032.25
$ws.Range("N134").Value = -68172.25

$ws.Range("H137").Value = 65275
$ws.Range("J137").Value = 65275
$ws.Range("L137").Value = 65275
$ws.Range("N137").Value = -75475

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 87324.46000000001
$ws.Range("I134").Value = 4225.1763
$ws.Range("J134").Value = 289137
$ws.Range("K134").Value = 12675.5289
$ws.Range("L134").Value = 867411
$ws.Range("M134").Value = -10140.5289
$ws.Range("N134").Value = -872481

$ws.Range("H135").Value = 63945.25
$ws.Range("J135").Value = 63945.25
$ws.Range("L135").Value = 63945.25
$ws.Range("N135").Value = -74085.25

$ws.Range("H140").Value = 79469.875
$ws.Range("J140").Value = 79469.875
$ws.Range("L140").Value = 79469.875
$ws.Range("N140").Value = -89829.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2310420.8
$ws.Range("I31").Value = 3002605.5
$ws.Range("J31").Value = 3138.889
$ws.Range("K31").Value = 3002605.5
$ws.Range("L31").Value = 3138.889
$ws.Range("M31").Value = -3002310.5
$ws.Range("N31").Value = -3728.889

$ws.Range("H34").Value = 2310420.8
$ws.Range("I34").Value = 3002605.5
$ws.Range("J34").Value = 3138.889
$ws.Range("K34").Value = 3002605.5
$ws.Range("L34").Value = 3138.889
$ws.Range("M34").Value = -3002403.5
$ws.Range("N34").Value = -3542.889

$ws.Range("H58").Value = 5667.04
$ws.Range("I58").Value = 1970.3572
$ws.Range("J58").Value = 10371.909
$ws.Range("K58").Value = 1970.3572
$ws.Range("L58").Value = 10371.909
$ws.Range("M58").Value = -1767.3572
$ws.Range("N58").Value = -10777.909

$ws.Range("H62").Value = 6382.5
$ws.Range("I62").Value = 6215.909
$ws.Range("J62").Value = 6993.3335
$ws.Range("K62").Value = 6215.909
$ws.Range("L62").Value = 6993.3335
$ws.Range("M62").Value = -5591.909
$ws.Range("N62").Value = -8241.333500000001

$ws.Range("H65").Value = 6382.5
$ws.Range("I65").Value = 6215.909
$ws.Range("J65").Value = 6993.3335
$ws.Range("K65").Value = 31079.545
$ws.Range("L65").Value = 34966.6675
$ws.Range("M65").Value = -27959.545
$ws.Range("N65").Value = -41206.6675

$ws.Range("H136").Value = 5667.04
$ws.Range("I136").Value = 1970.3572
$ws.Range("J136").Value = 10371.909
$ws.Range("K136").Value = 5911.071599999999
$ws.Range("L136").Value = 31115.727
$ws.Range("M136").Value = -3361.071599999999
$ws.Range("N136").Value = -36215.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 6000
$ws.Range("J104").Value = 6000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -23242

$ws.Range("H113").Value = 550.125
$ws.Range("I113").Value = 970.4286
$ws.Range("J113").Value = 478.36584
$ws.Range("K113").Value = 2911.2858
$ws.Range("L113").Value = 1435.09752
$ws.Range("M113").Value = -741.2857999999997
$ws.Range("N113").Value = -5775.09752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 55333.332
$ws.Range("J133").Value = 55333.332
$ws.Range("L133").Value = 55333.332
$ws.Range("N133").Value = -65453.332

$ws.Range("H135").Value = 44853.332
$ws.Range("J135").Value = 44853.332
$ws.Range("L135").Value = 44853.332
$ws.Range("N135").Value = -54993.332

$ws.Range("H138").Value = 54750
$ws.Range("J138").Value = 54750
$ws.Range("L138").Value = 54750
$ws.Range("N138").Value = -65030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 742.5238000000001
$ws.Range("I22").Value = 408.5
$ws.Range("J22").Value = 1187.8889
$ws.Range("K22").Value = 408.5
$ws.Range("L22").Value = 1187.8889
$ws.Range("M22").Value = -113.5
$ws.Range("N22").Value = -1777.8889

$ws.Range("H27").Value = 742.5238000000001
$ws.Range("I27").Value = 408.5
$ws.Range("J27").Value = 1187.8889
$ws.Range("K27").Value = 408.5
$ws.Range("L27").Value = 1187.8889
$ws.Range("M27").Value = -301.5
$ws.Range("N27").Value = -1401.8889

$ws.Range("H133").Value = 43000
$ws.Range("J133").Value = 43000
$ws.Range("L133").Value = 43000
$ws.Range("N133").Value = -48060

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 3796.238
$ws.Range("I136").Value = 2183.818
$ws.Range("J136").Value = 5569.9
$ws.Range("K136").Value = 6551.454000000001
$ws.Range("L136").Value = 16709.7
$ws.Range("M136").Value = -4001.454000000001
$ws.Range("N136").Value = -21809.7

$ws.Range("H141").Value = 95800
$ws.Range("J141").Value = 95800
$ws.Range("L141").Value = 95800
$ws.Range("N141").Value = -106160

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 28331.666
$ws.Range("J110").Value = 28331.666
$ws.Range("L110").Value = 28331.666
$ws.Range("N110").Value = -36511.666

$ws.Range("H141").Value = 94000
$ws.Range("J141").Value = 94000
$ws.Range("L141").Value = 94000
$ws.Range("N141").Value = -104360
